$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.077473474454563984,
    -0.0099999998471815843,
    -0.018082176307787634,
    0.28399768752469967,
    -0.0059999998522961562,
    -0.005999999846913795,
    -0.019999999820628389,
    -0.019999999818497649,
    -0.0059999998416042644,
    -0.0081760631161955644,
    -0.0044999998417303289,
    -0.0059999998384370201,
    -0.0059999998358462037,
    -0.011999999824148233,
    -0.0059999998345743322,
    -0.005999999834060965,
    -0.0059999998333930549,
    -0.0089999998277496829,
    -0.0089999998488909938,
    -0.0089999998476333332,
    -0.0089999998474654674,
    -0.0089999998473393461,
    -0.0089999998472363174,
    -0.041999999784455788,
    -0.041999999783368658,
    -0.0059999998464874693,
    -0.0059999998459616677,
    -0.005999999843653292,
    0.011443851331305765,
    -0.01999999981547429,
    -0.021237953633107765,
    -0.020999999811714964,
    -0.0059999998390924958
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
